$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.920.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.888.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7723"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3099"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08549"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.45%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7641"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.900.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.346"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.158"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.057.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.268.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007805"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.962"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1635"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.323"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.034"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.437"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.502"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.107"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05434"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.241"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7474"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.004"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.696"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01964"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.783"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4461"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.110.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.92%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.092"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8487"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.155.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.871"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.600"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.991"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.04%  "
